$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160; existing rows 160:183 shift down to 161:184
# (matches the diff's dimension change from A1:R183 to A1:R184, with every
# existing data row moving down by one and a brand-new weekly record landing
# on row 160).
$ws.Rows.Item(160).Insert()

$ws.Range("A160").Value = 8
$ws.Range("B160").Value = "Terminal La Palmera de La Serena"
$ws.Range("C160").Value = "Coquimbo"
$ws.Range("D160").Value = 44522
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 100112012
$ws.Range("G160").Value = "Espinaca"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 2000
$ws.Range("K160").Value = 400
$ws.Range("L160").Value = 500
$ws.Range("M160").Value = 450
$ws.Range("N160").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O160").Value = "Provincia del Elquí"
$ws.Range("P160").Value = 900
$ws.Range("Q160").Value = 0.5
$ws.Range("R160").Value = "Hortaliza"
